# "ajout dans tableau de bord" - add a new contribution row (row 10) to
# the dashboard table on the single worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of the previous data row (row 9) down onto the new
# row 10 so the date/name cells keep the same number format & style.
$ws.Range("B9:E9").Copy()
$ws.Range("B10").PasteSpecial(-4122)  # xlPasteFormats

# B10: Date -> 06/01/2021 (same day as row 9, serial 44202)
$ws.Range("B10").Value2 = 44202

# C10: Nom -> François
$ws.Range("C10").Value = "François"

# D10: Fonctionnalité -> n°18
$ws.Range("D10").Value = 18

# E10: Contribution -> new note
$ws.Range("E10").Value = "Finalisation état authentifié"

$wb.Save()
